$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("vasenate")

# New candidate rows to append after the existing data (rows 2-62 -> new rows 63-65)
# Fill in id/first/last names first (matches the order new shared strings were created)
$ws.Range("A63").Value = 2161
$ws.Range("B63").Value = "Wayne"
$ws.Range("D63").Value = "Powell"

$ws.Range("A64").Value = 2171
$ws.Range("B64").Value = "Amy"
$ws.Range("D64").Value = "Laufer"

$ws.Range("A65").Value = 2181
$ws.Range("B65").Value = "Geary"
$ws.Range("D65").Value = "Higgins"

# Party id / office / district / in EC? / election id
$ws.Range("F63").Value = "dem"
$ws.Range("G63").Value = "vasenate"
$ws.Range("H63").Value = 11
$ws.Range("I63").Value = "yes"
$ws.Range("J63").Value = 1871

$ws.Range("F64").Value = "dem"
$ws.Range("G64").Value = "vasenate"
$ws.Range("H64").Value = 17
$ws.Range("I64").Value = "yes"
$ws.Range("J64").Value = 1581

$ws.Range("F65").Value = "rep"
$ws.Range("G65").Value = "vasenate"
$ws.Range("H65").Value = 13
$ws.Range("I65").Value = "yes"
$ws.Range("J65").Value = 1781

# Websites last: plain text for Powell, hyperlink for Laufer
$ws.Range("E63").Value = "https://www.powellsenateva.com/"
$ws.Hyperlinks.Add($ws.Range("E64"), "https://www.lauferforvirginia.com/")

# Remove the autofilter that previously covered A1:J62
if ($ws.AutoFilterMode) {
    $ws.AutoFilterMode = $false
}

# Update selection / view state and make vasenate the active sheet/tab
$ws.Activate()
[void]$ws.Range("I66").Select()
$excel.ActiveWindow.ScrollRow = 36
